$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.250.36'
$ws.Range('E2').Value = '  +1.81%  '
$ws.Range('D3').Value = '2.425.93'
$ws.Range('E3').Value = '  +2.14%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '561.26'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.37%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.34'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.18%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('E8').Value = '  +2.39%  '
$ws.Range('D9').Value = '2.424.13'
$ws.Range('E9').Value = '  +1.95%  '
$ws.Range('E10').Value = '  +1.56%  '
$ws.Range('E11').Value = '  -1.55%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.42'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.76%  '
$ws.Range('E13').Value = '  +1.69%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.20'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.60%  '
$ws.Range('E15').Value = '  +6.04%  '
$ws.Range('D16').Value = '2.853.79'
$ws.Range('E16').Value = '  +2.36%  '
$ws.Range('D17').Value = '62.129.77'
$ws.Range('E17').Value = '  +1.77%  '
$ws.Range('D18').Value = '2.424.67'
$ws.Range('E18').Value = '  +2.00%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.23'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.93%  '
$ws.Range('E20').Value = '  +1.45%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '324.85'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.42%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.81'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.61%  '
$ws.Range('E23').Value = '  +0.10%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '65.56'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.92%  '
$ws.Range('E25').Value = '  +1.53%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.96'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +5.54%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '594.25'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +17.03%  '
$ws.Range('D28').Value = '2.540.60'
$ws.Range('E28').Value = '  +1.98%  '
$ws.Range('E29').Value = '  +0.16%  '
$ws.Range('D30').Value = '0.0₃0943'
$ws.Range('E30').Value = '  +5.89%  '
$ws.Range('E31').Value = '  +2.03%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.46'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +6.13%  '
$ws.Range('E33').Value = '  -0.11%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.89'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.45%  '
$ws.Range('E35').Value = '  +2.79%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.75'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +6.14%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.00'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.02%  '
$ws.Range('E38').Value = '  +2.95%  '
$ws.Range('E39').Value = '  +1.99%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '153.39'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +4.28%  '
$ws.Range('E41').Value = '  +1.64%  '
$ws.Range('E42').Value = '  -3.93%  '
$ws.Range('E43').Value = '  +0.08%  '
$ws.Range('E44').Value = '  +11.35%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '151.41'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.07%  '
$ws.Range('E46').Value = '  +1.99%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0542'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +4.13%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '20.37'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +5.87%  '
$ws.Range('E49').Value = '  +2.93%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0923'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.86%  '
$ws.Range('E51').Value = '  +2.89%  '
